# Update countries & provincias Spain
# Applies the diff to the "Pais" sheet of paises.xlsx:
#  - swap the "Japon"/"Luxemburgo" rows (incl. their stats) at rows 34-35
#  - swap the "Moldavia"/"Kazajistan" rows (incl. their stats) at rows 78-79
#  - refresh the "Datos actualizados..." timestamp string in A1
#  - update the various per-country statistic cells (columns B:H) with the
#    latest counts

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (A1) -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 07:50"

# --- Row 34 / 35: swap Japon <-> Luxemburgo (country names + their stats) -
$ws.Range("A34").Value = "Luxemburgo"
$ws.Range("B34").Value = 2178
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 80
$ws.Range("E34").Value = 2075
$ws.Range("F34").Value = 31
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 23

$ws.Range("A35").Value = "Japon"
$ws.Range("B35").Value = 2178
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 424
$ws.Range("E35").Value = 1697
$ws.Range("F35").Value = 69
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 57

# --- Row 78 / 79: swap Moldavia <-> Kazajistan (country names + stats) ----
$ws.Range("A78").Value = "Kazajistan"
$ws.Range("B78").Value = 363
$ws.Range("C78").Value = 20
$ws.Range("D78").Value = 24
$ws.Range("E78").Value = 336
$ws.Range("F78").Value = 6
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = 3

$ws.Range("A79").Value = "Moldavia"
$ws.Range("B79").Value = 353
$ws.Range("C79").Value = 0
$ws.Range("D79").Value = 18
$ws.Range("E79").Value = 331
$ws.Range("F79").Value = 44
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 4

# --- Remaining statistic refreshes -----------------------------------------

# Estados Unidos (row 4)
$ws.Range("E4").Value = 177272
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 4055

# China (row 7)
$ws.Range("B7").Value = 81554
$ws.Range("C7").Value = 36
$ws.Range("D7").Value = 76238
$ws.Range("E7").Value = 2004
$ws.Range("F7").Value = 466
$ws.Range("G7").Value = 7
$ws.Range("H7").Value = 3312

# Australia (row 22)
$ws.Range("B22").Value = 4860
$ws.Range("C22").Value = 97
$ws.Range("E22").Value = 4494
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 21

# Pakistan (row 37)
$ws.Range("B37").Value = 2039
$ws.Range("C37").Value = 101
$ws.Range("D37").Value = 82
$ws.Range("E37").Value = 1931

# Hungria (row 70)
$ws.Range("B70").Value = 525
$ws.Range("C70").Value = 33
$ws.Range("D70").Value = 40
$ws.Range("E70").Value = 469

# Bulgaria (row 73)
$ws.Range("B73").Value = 412
$ws.Range("C73").Value = 13
$ws.Range("E73").Value = 387

# Uzbekistan (row 99)
$ws.Range("B99").Value = 173
$ws.Range("C99").Value = 1
$ws.Range("D99").Value = 8

# Nigeria (row 107)
$ws.Range("B107").Value = 139
$ws.Range("C107").Value = 4
$ws.Range("D107").Value = 9
$ws.Range("E107").Value = 128

# Trinidad yTobago (row 120)
$ws.Range("B120").Value = 89
$ws.Range("C120").Value = 2
$ws.Range("E120").Value = 84
